$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("FB")
$ws4 = $wb.Worksheets.Item("FB rieng")

# The edit swaps the B:C login-data block between sheet "FB" (rows 2-4) and
# sheet "FB rieng" (rows 8-13). "FB rieng" has twice as many rows of data as
# "FB" currently holds, so the destination ranges are not the same size.
# Move through a temporary staging area on "FB" so both halves relocate
# cleanly in one pass without a same-sheet/overlapping-size ambiguity.

# Step 1: stash "FB" B2:C4 (current 3 rows of creds) in a scratch area.
$ws2.Range("B2:C4").Cut($ws2.Range("H2:I4"))

# Step 2: move "FB rieng" B8:C13 (6 rows of creds) up into "FB" B2:C7.
$ws4.Range("B8:C13").Cut($ws2.Range("B2:C7"))

# Step 3: move the stashed "FB" creds down into "FB rieng" B8:C10.
$ws2.Range("H2:I4").Cut($ws4.Range("B8:C10"))

# Step 4: drop the now-empty scratch cells so they don't linger in the sheet.
$ws2.Range("H2:I4").Clear()

# Restore the recorded cursor/selection position on each affected sheet.
$ws4.Activate()
$ws4.Range("B8:C10").Select()

$ws2.Activate()
$ws2.Range("B6").Select()
